# "No474. Ones and Zeroes finished"
# Append a new tracker row (row 28) for LeetCode 474. Ones and Zeroes,
# mirroring the structure/format of the existing rows 25-27, wire up its
# hyperlink, and move the view/selection the way the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fill in the new row's values --------------------------------------
$ws.Range("A28").Value = "474. Ones and Zeroes"
$ws.Range("B28").Value = "Medium"
$ws.Range("C28").Value = "https://leetcode.com/problems/ones-and-zeroes/"
$ws.Range("D28").Value = 44506
$ws.Range("E28").Value = "0-1背包问题"
$ws.Range("F28").Value = "背包体积约束条件有两个，dp数组多一个维度"
$ws.Range("G28").Value = "未复习"
$ws.Range("H28").Value = "⭕"

# --- 2. Hyperlink the problem URL in C28 -----------------------------------
$ws.Hyperlinks.Add($ws.Range("C28"), "https://leetcode.com/problems/ones-and-zeroes/")

# --- 3. Re-apply row 26's formatting (same style pattern the new row needs)
#        AFTER adding the hyperlink, since Hyperlinks.Add stamps its own
#        style onto the anchor cell.
$ws.Range("A26:H26").Copy()
$ws.Range("A28:H28").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- 4. Row height to match the other short (single-line) rows -------------
$ws.Rows.Item(28).RowHeight = 28

# --- 5. Leave the selection/scroll where the author left it ----------------
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("I25").Select()
